$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.927.83"
$ws.Range("E2").Value = "  +0.77%  "

$ws.Range("D3").Value = "1.633.44"
$ws.Range("E3").Value = "  +1.95%  "

$ws.Range("E4").Value = "  +0.30%  "

$ws.Range("D5").Value = "'214.42"
$ws.Range("E5").Value = "  +0.88%  "

$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("E7").Value = "  +0.24%  "

$ws.Range("D8").Value = "'28.55"
$ws.Range("E8").Value = "  +1.23%  "

$ws.Range("E9").Value = "  +1.43%  "

$ws.Range("E10").Value = "  +0.74%  "

$ws.Range("D11").Value = "'0.0907"
$ws.Range("E11").Value = "  -0.02%  "

$ws.Range("D12").Value = "1.868.65"
$ws.Range("E12").Value = "  +2.04%  "

$ws.Range("D13").Value = "1.636.05"
$ws.Range("E13").Value = "  +2.16%  "

$ws.Range("D14").Value = "'0.563"
$ws.Range("E14").Value = "  +1.97%  "

$ws.Range("D15").Value = "'9.35"
$ws.Range("E15").Value = "  +17.66%  "

$ws.Range("D16").Value = "29.953.39"
$ws.Range("E16").Value = "  +0.80%  "

$ws.Range("D17").Value = "'3.85"
$ws.Range("E17").Value = "  +2.30%  "

$ws.Range("D18").Value = "'64.07"
$ws.Range("E18").Value = "  +0.08%  "

$ws.Range("D19").Value = "'242.23"
$ws.Range("E19").Value = "  +0.14%  "

$ws.Range("E20").Value = "  +0.52%  "

$ws.Range("D21").Value = "'1.00"
$ws.Range("E21").Value = "  +0.14%  "

$ws.Range("D22").Value = "'9.83"
$ws.Range("E22").Value = "  +4.28%  "

$ws.Range("D23").Value = "'4.14"
$ws.Range("E23").Value = "  +2.61%  "

$ws.Range("E24").Value = "  +2.69%  "

$ws.Range("D25").Value = "'157.74"
$ws.Range("E25").Value = "  +1.72%  "

$ws.Range("D26").Value = "'15.52"
$ws.Range("E26").Value = "  +0.39%  "

$ws.Range("E27").Value = "  +0.56%  "

$ws.Range("D28").Value = "'6.62"
$ws.Range("E28").Value = "  +2.60%  "

$ws.Range("E29").Value = "  +0.22%  "

$ws.Range("E30").Value = "  +1.75%  "

$ws.Range("E31").Value = "  +4.47%  "

$ws.Range("E32").Value = "  +4.20%  "

$ws.Range("E33").Value = "  -0.39%  "

$ws.Range("D34").Value = "1.425.34"
$ws.Range("E34").Value = "  +0.33%  "

$ws.Range("E35").Value = "  +5.08%  "

$ws.Range("E36").Value = "  -0.26%  "

$ws.Range("E37").Value = "  -3.04%  "

$ws.Range("E38").Value = "  -0.26%  "

$ws.Range("E39").Value = "  +0.40%  "

$ws.Range("D40").Value = "'76.09"
$ws.Range("E40").Value = "  +13.10%  "

$ws.Range("E41").Value = "  +1.50%  "

$ws.Range("E42").Value = "  +2.96%  "

$ws.Range("D43").Value = "'0.830"
$ws.Range("E43").Value = "  +1.72%  "

$ws.Range("D44").Value = "'0.0491"
$ws.Range("E44").Value = "  -0.58%  "

$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'1.02"
$ws.Range("E45").Value = "  +2.75%  "

$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").Value = "'52.95"
$ws.Range("E47").Value = "  -4.66%  "

$ws.Range("E48").Value = "  -0.20%  "

$ws.Range("D49").Value = "1.775.29"
$ws.Range("E49").Value = "  +2.09%  "

$ws.Range("D50").Value = "'90.59"
$ws.Range("E50").Value = "  +4.54%  "

$ws.Range("D51").Value = "0.0₆0113"
$ws.Range("E51").Value = "  +8.46%  "
